# Extent Report Scripts added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC02 / TC04 Trigger column: N -> Y
$ws.Range("C8").Value = "Y"
$ws.Range("C20").Value = "Y"

# Update verifyURL "Value" cells and turn them into hyperlinks.
# Order matches the target relationship ids (F28 -> rId1, F22 -> rId2, F10 -> rId3).
$ws.Hyperlinks.Add($ws.Range("F28"), "https://github.com/login")

$ws.Range("F22").Value = "https://github.com/login04"
$ws.Hyperlinks.Add($ws.Range("F22"), "https://github.com/login04")

$ws.Range("F10").Value = "https://github.com/login04"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://github.com/login04")

# Move the sheet selection to H4
$ws.Range("H4").Select()
